# Publish terminology IG 2.0.0:
#  - Metadata: bump Version / Date
#  - Split the old "Concepts" sheet into a new "Properties" sheet (FHIR
#    CodeSystem .property definitions) plus a "Concepts" sheet that keeps
#    the original concept table (Level/Code/Display/Definition).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: Version 1.8.1 -> 1.1.1, Date 2023-10-31 -> 2025-09-22
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value2 = "1.1.1"
# Force the date cell to stay plain text (matches the source workbook,
# where it is a shared string, not a date serial number).
$meta.Range("B8").NumberFormat = "@"
$meta.Range("B8").Value2 = "2025-09-22"

# ---------------------------------------------------------------------
# 2. Insert a new sheet right after "Concepts" that will become the new
#    "Concepts" sheet; copy the existing concepts table (with formatting)
#    into it, then rename the sheets so the old "Concepts" sheet becomes
#    "Properties".
# ---------------------------------------------------------------------
$oldConcepts = $wb.Worksheets.Item("Concepts")

$newConcepts = $wb.Worksheets.Add($null, $oldConcepts)
$oldConcepts.Range("A1:D6").Copy($newConcepts.Range("A1"))

$oldConcepts.Name = "Properties"
$newConcepts.Name = "Concepts"

# ---------------------------------------------------------------------
# 3. Replace the "Properties" sheet content with the CodeSystem property
#    definitions (Code / Uri / Description / Type).
# ---------------------------------------------------------------------
$props = $wb.Worksheets.Item("Properties")

# Drop rows 4:6 inherited from the old Concepts data - Properties only
# needs a header row + 2 data rows.
$props.Range("A4:D6").Delete()

$props.Range("A1").Value2 = "Code"
$props.Range("B1").Value2 = "Uri"
$props.Range("C1").Value2 = "Description"
$props.Range("D1").Value2 = "Type"

$props.Range("A2").Value2 = "status"
$props.Range("B2").Value2 = "http://hl7.org/fhir/concept-properties#status"
$props.Range("C2").Value2 = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$props.Range("D2").Value2 = "code"

$props.Range("A3").Value2 = "effectiveDate"
$props.Range("B3").Value2 = "http://hl7.org/fhir/concept-properties#effectiveDate"
$props.Range("C3").Value2 = "The date at which the concept status was last changed."
$props.Range("D3").Value2 = "dateTime"

# Keep the originally-selected tab ("Metadata") active.
$meta.Activate()
